# Update Overview Decks for August 2025 (#62)
#
# 1) Re-point every "Member Benefits" style table (6 of them, scattered
#    across the deck) at the new table style GUID.
# 2) Swap the two unused/"spare" theme colour palettes that ship in the
#    package (the ones behind the "Simple Light" design / notes master)
#    so the RGB values that used to live under the Default palette now
#    live under Simple Light, and vice-versa.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Table style id swap
# ---------------------------------------------------------------------
$oldTableStyle = "{98510C90-4868-4D3F-B63A-26FF2F9ABE26}"
$newTableStyle = "{22567273-7B7A-4A14-B2B7-553CC2E4EEBA}"

for ($n = 1; $n -le $p.Slides.Count; $n++) {
    $slide = $p.Slides.Item($n)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ("$($table.Style)" -eq $oldTableStyle) {
                $table.ApplyStyle($newTableStyle, $true)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Theme colour scheme swap (Simple Light <-> Default)
# ---------------------------------------------------------------------
# Colours that used to be the "Default" scheme -- becomes Simple Light's.
$defaultScheme = @(0, 16777215, 5800213, 15987699, 13077765, 3322960, 1791725, 61421, 15059748, 7529828, 13369378, 9116245)

$master = $p.Designs.Item(1).SlideMaster
$scheme = $master.ColorScheme
for ($idx = 1; $idx -le 12; $idx++) {
    $scheme.Colors($idx).RGB = $defaultScheme[$idx - 1]
}
